$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.240.06"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.596.85"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'315.80"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'98.14"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "'35.93"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'7.55"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.990.66"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "2.486.96"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "43.326.66"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'69.62"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "'255.00"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("D26").Value = "'27.32"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'41.05"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'10.33"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "'5.88"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'156.01"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'3.47"
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("D35").Value = "'0.0811"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("D36").Value = "'2.71"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("D37").Value = "'18.85"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +8.35%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "'23.11"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("E42").Value = "  +6.56%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D45").Value = "'3.25"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "2.013.39"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'8.99"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'83.64"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.780.42"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.196"
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'74.93"
$ws.Range("E51").Value = "  +1.96%  "
